# Add a new "2022" data column (K) to the indicators table on sheet1,
# mirroring the formatting already used by the neighbouring "2020" column (I)
# and filling in the new year's figures (leaving the section-header rows
# blank, same as the other year columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone the formatting of column I (rows 4-46, header + data rows) into
#    column K so the new column visually matches the rest of the table
#    (number format, font, alignment, etc.) before we pour in the values.
$ws.Range("I4:I46").Copy($ws.Range("K4:K46"))

# 2) Header for the new year.
$ws.Range("K4").Value = 2022

# 3) Data values for the new "2022" column. Section-title rows (6, 9, 22,
#    25) intentionally stay blank, matching the other year columns.
$ws.Range("K5").Value = 24.2
$ws.Range("K7").Value = 25.5
$ws.Range("K8").Value = 22.3
$ws.Range("K10").Value = 18
$ws.Range("K11").Value = 18.899999999999999
$ws.Range("K12").Value = 21.3
$ws.Range("K13").Value = 30.2
$ws.Range("K14").Value = 31.7
$ws.Range("K15").Value = 34.1
$ws.Range("K16").Value = 25.8
$ws.Range("K17").Value = 20
$ws.Range("K18").Value = 12.1
$ws.Range("K19").Value = 10.3
$ws.Range("K20").Value = 15.1
$ws.Range("K21").Value = 12.1
$ws.Range("K23").Value = 25.9
$ws.Range("K24").Value = 23.2
$ws.Range("K26").Value = 25.9
$ws.Range("K27").Value = 48.3
$ws.Range("K28").Value = 24.3
$ws.Range("K29").Value = 28.1
$ws.Range("K30").Value = 25.8
$ws.Range("K31").Value = 27.1
$ws.Range("K32").Value = 20.7
$ws.Range("K33").Value = 24.3
$ws.Range("K34").Value = 19.399999999999999
$ws.Range("K35").Value = 7.5
$ws.Range("K36").Value = 11.4
$ws.Range("K37").Value = 36.5
$ws.Range("K38").Value = 17.8
$ws.Range("K39").Value = 20.3
$ws.Range("K40").Value = 20.5
$ws.Range("K41").Value = 32.200000000000003
$ws.Range("K42").Value = 23.2
$ws.Range("K43").Value = 23.8
$ws.Range("K44").Value = 21
$ws.Range("K45").Value = 18
$ws.Range("K46").Value = 3.2

# 4) Move the active selection like the author left it.
$ws.Range("L12").Select() | Out-Null
